# Se crea la reimpresion de comprobantes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# P3: comprobante number label (was "N°:267" -> now "N°:15")
$ws.Range("P3").Value = "N°:15"

# E8: cantidad / multiplicador (1 -> 3)
$ws.Range("E8").Value = 3

# J11: "La cantidad de" text amount (was "UN PESO " -> now new total in words)
$ws.Range("J11").Value = "TRESCIENTOS SEIS MIL SETECIENTOS NOVENTA Y SEIS PESOS "

# P15: cantidad efectivo/transferencia (1 -> 3)
$ws.Range("P15").Value = 3

# K18: observaciones (empty -> "Sin observaciones")
$ws.Range("K18").Value = "Sin observaciones"
